# Apply the text replacements required by the diff.
# All "old" strings are unique within the document and none of the
# "new" strings collide with any other "old" string, so a simple
# sequential Find/Replace across the whole document content is safe.

$d = $word.ActiveDocument

$replacements = @(
    @("2025-06-10 Tuesday", "2025-06-11 Wednesday"),
    @("93×91=", "35×94="),
    @("46×39=", "98×88="),
    @("31×50=", "61×91="),
    @("88×69=", "92×26="),
    @("41×26=", "59×49="),
    @("27×18=", "19×26="),
    @("25×54=", "60×40="),
    @("25×43=", "35×54="),
    @("18×84=", "80×90="),
    @("96×30=", "39×85="),
    @("89×94=", "62×39="),
    @("80×32=", "31×16="),
    @("17×19=", "98×57="),
    @("55×57=", "62×70="),
    @("58×52=", "67×40="),
    @("14×79=", "93×94="),
    @("46×42=", "17×22="),
    @("25×73=", "23×34="),
    @("11×23=", "34×19="),
    @("26×82=", "54×30="),
    @("62×24=", "76×86="),
    @("12×66=", "43×15="),
    @("96×80=", "86×79="),
    @("76×98=", "36×82="),
    @("86×99=", "74×88=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
